# Update "想去人数" (want-to-go count) figures across the workbook's sheets.
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 103
$wsExhibit.Range("F3").Value = 12083
$wsExhibit.Range("F4").Value = 45
$wsExhibit.Range("F8").Value = 11972
$wsExhibit.Range("F10").Value = 1183
$wsExhibit.Range("F13").Value = 1798
$wsExhibit.Range("F14").Value = 5924
$wsExhibit.Range("F15").Value = 131
$wsExhibit.Range("F16").Value = 3557
$wsExhibit.Range("F17").Value = 203

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 9

# Sheet "全部类型" (All types, merged view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 103
$wsAll.Range("F5").Value = 12083
$wsAll.Range("F6").Value = 45
$wsAll.Range("F8").Value = 9
$wsAll.Range("F11").Value = 11972
$wsAll.Range("F13").Value = 1183
$wsAll.Range("F16").Value = 1798
$wsAll.Range("F18").Value = 5924
$wsAll.Range("F19").Value = 131
$wsAll.Range("F20").Value = 3557
$wsAll.Range("F21").Value = 203
